# Update the cryptos list (Price / Volume(1h) columns) with the latest
# scraped figures. Price values that look like plain decimal numbers are
# prefixed with a leading apostrophe so Excel stores them as text (matching
# the source data, which is text even when it happens to parse as a
# number) instead of silently converting them to numeric cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.434.22"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "3.945.63"
$ws.Range("E3").Value = "  +4.40%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'469.84"
$ws.Range("E5").Value = "  +8.18%  "
$ws.Range("D6").Value = "'146.35"
$ws.Range("E6").Value = "  +4.56%  "
$ws.Range("D7").Value = "'0.625"
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.732"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  +11.31%  "
$ws.Range("D11").Value = "'0.0000356"
$ws.Range("E11").Value = "  +13.36%  "
$ws.Range("D12").Value = "'43.38"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "4.575.16"
$ws.Range("E13").Value = "  +3.91%  "
$ws.Range("D14").Value = "'10.38"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "'15.17"
$ws.Range("E15").Value = "  +2.28%  "
$ws.Range("D16").Value = "3.939.90"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("D18").Value = "'19.88"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'1.15"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "67.639.75"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "'435.14"
$ws.Range("E21").Value = "  +4.06%  "
$ws.Range("D22").Value = "'3.39"
$ws.Range("E22").Value = "  +4.33%  "
$ws.Range("D23").Value = "'14.39"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "'87.45"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").Value = "'3.62"
$ws.Range("E25").Value = "  +7.21%  "
$ws.Range("D26").Value = "'38.64"
$ws.Range("E26").Value = "  +3.99%  "
$ws.Range("D27").Value = "'10.19"
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").Value = "'721.07"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").Value = "'0.132"
$ws.Range("E30").Value = "  -0.62%  "
$ws.Range("D31").Value = "'13.45"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("D33").Value = "'42.23"
$ws.Range("E33").Value = "  -2.78%  "
$ws.Range("D34").Value = "0.0₃0856"
$ws.Range("E34").Value = "  +26.34%  "
$ws.Range("D35").Value = "'57.94"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("D36").Value = "'0.150"
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'5.32"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("E40").Value = "  +4.85%  "
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("E42").Value = "  +6.73%  "
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'0.335"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("E45").Value = "  +6.61%  "
$ws.Range("D46").Value = "'2.21"
$ws.Range("E46").Value = "  +6.39%  "
$ws.Range("E47").Value = "  -5.60%  "
$ws.Range("D48").Value = "'147.91"
$ws.Range("E48").Value = "  +3.79%  "
$ws.Range("D49").Value = "'3.17"
$ws.Range("E49").Value = "  -4.90%  "
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").Value = "'25.44"
$ws.Range("E51").Value = "  +2.81%  "
